$wb = $excel.ActiveWorkbook

# --- Typography sheet: fill in Wildcard Characters / Wildcard Ranges for the
#     first (Default/ADC) row, restricting it to digits only.
$ws1 = $wb.Worksheets.Item("Typography")
$ws1.Range("G4").Value = "0-9"
$ws1.Range("H4").Value = "0-9"

# --- Translation sheet: add the three new ADC-related translation rows
#     (TEXT ID / TYPOGRAPHY NAME / ALIGNMENT / DIRECTION / GB).
$ws2 = $wb.Worksheets.Item("Translation")

$ws2.Range("B4").Value = "SingleUseId1"
$ws2.Range("C4").Value = "Default"
$ws2.Range("D4").Value = "Left"
$ws2.Range("E4").Value = "LTR"
$ws2.Range("F4").Value = "ADC value = <value>"

$ws2.Range("B5").Value = "SingleUseId2"
$ws2.Range("C5").Value = "Default"
$ws2.Range("D5").Value = "Left"
$ws2.Range("E5").Value = "LTR"
# "10" looks numeric, so a plain .Value assignment would store it as a
# number; force it to land as text (matching the source workbook) via a
# formula-to-value round trip instead of touching NumberFormat (which would
# leave a stray cell style behind).
$f5 = $ws2.Range("F5")
$f5.Formula = "=""10"""
$f5.Copy()
$f5.PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws2.Range("B6").Value = "SingleUseId3"
$ws2.Range("C6").Value = "Default"
$ws2.Range("D6").Value = "Center"
$ws2.Range("E6").Value = "LTR"
$ws2.Range("F6").Value = "New ADC"
